$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-10 (NATMI re-run results), columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T.
# F and L columns are unchanged.
$updates = @{
    "2" = @{ "E" = 3; "G" = 21.84976866666667; "H" = 65.549306; "I" = 0.05020018890879543; "J" = 0.05020018890879543; "K" = 3; "M" = 0.8749903333333333; "N" = 2.624971; "O" = 0.2670516933349977; "P" = 0.2670516933349977; "Q" = 19.11833636890289; "R" = 172.065027320126; "S" = 0.01340604545383059; "T" = 0.01340604545383059 }
    "3" = @{ "E" = 3; "G" = 21.84976866666667; "H" = 65.549306; "I" = 0.05020018890879543; "J" = 0.05020018890879543; "K" = 3; "M" = 0.8147036666666666; "N" = 2.444111; "O" = 0.2486518827250642; "P" = 0.2486518827250642; "Q" = 17.80108664855178; "R" = 160.209779836966; "S" = 0.01248237148532587; "T" = 0.01248237148532587 }
    "4" = @{ "E" = 3; "G" = 21.84976866666667; "H" = 65.549306; "I" = 0.05020018890879543; "J" = 0.05020018890879543; "K" = 3; "M" = 1.586789; "N" = 4.760367; "O" = 0.484296423939938; "P" = 0.484296423939938; "Q" = 34.67097257281133; "R" = 312.038753155302; "S" = 0.02431177196963896; "T" = 0.02431177196963896 }
    "5" = @{ "E" = 3; "G" = 385.0524703333334; "H" = 1155.157411; "I" = 0.8846641374295412; "J" = 0.8846641374295412; "K" = 3; "M" = 0.8749903333333333; "N" = 2.624971; "O" = 0.2670516933349977; "P" = 0.2670516933349977; "Q" = 336.9171893677868; "R" = 3032.254704310081; "S" = 0.2362510559333041; "T" = 0.2362510559333041 }
    "6" = @{ "E" = 3; "G" = 385.0524703333334; "H" = 1155.157411; "I" = 0.8846641374295412; "J" = 0.8846641374295412; "K" = 3; "M" = 0.8147036666666666; "N" = 2.444111; "O" = 0.2486518827250642; "P" = 0.2486518827250642; "Q" = 313.7036594396246; "R" = 2823.332934956621; "S" = 0.2199734033512004; "T" = 0.2199734033512004 }
    "7" = @{ "E" = 3; "G" = 385.0524703333334; "H" = 1155.157411; "I" = 0.8846641374295412; "J" = 0.8846641374295412; "K" = 3; "M" = 1.586789; "N" = 4.760367; "O" = 0.484296423939938; "P" = 0.484296423939938; "Q" = 610.9970243477596; "R" = 5498.973219129837; "S" = 0.4284396781450366; "T" = 0.4284396781450366 }
    "8" = @{ "E" = 3; "G" = 28.350479; "H" = 85.05143699999999; "I" = 0.06513567366166337; "J" = 0.06513567366166337; "K" = 3; "M" = 0.8749903333333333; "N" = 2.624971; "O" = 0.2670516933349977; "P" = 0.2670516933349977; "Q" = 24.80639507036966; "R" = 223.257555633327; "S" = 0.01739459194786301; "T" = 0.01739459194786301 }
    "9" = @{ "E" = 3; "G" = 28.350479; "H" = 85.05143699999999; "I" = 0.06513567366166337; "J" = 0.06513567366166337; "K" = 3; "M" = 0.8147036666666666; "N" = 2.444111; "O" = 0.2486518827250642; "P" = 0.2486518827250642; "Q" = 23.09723919305633; "R" = 207.875152737507; "S" = 0.01619610788853798; "T" = 0.01619610788853797 }
    "10" = @{ "E" = 3; "G" = 28.350479; "H" = 85.05143699999999; "I" = 0.06513567366166337; "J" = 0.06513567366166337; "K" = 3; "M" = 1.586789; "N" = 4.760367; "O" = 0.484296423939938; "P" = 0.484296423939938; "Q" = 44.98622822193099; "R" = 404.8760539973789; "S" = 0.03154497382526238; "T" = 0.03154497382526238 }
}

foreach ($r in $updates.Keys) {
    $rowVals = $updates[$r]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$r").Value = $rowVals[$col]
    }
}